$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.975.11'
$ws.Range('E2').Value = '  +1.60%  '

# Row 3
$ws.Range('D3').Value = '3.149.13'
$ws.Range('E3').Value = '  +2.75%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.75%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.49%  '

# Row 7
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$ws.Range('D8').Value = '3.151.89'
$ws.Range('E8').Value = '  +2.89%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.70%  '

# Row 10
$ws.Range('E10').Value = '  +3.84%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.11'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.77%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.498'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.34%  '

# Row 13
$ws.Range('E13').Value = '  +12.92%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.57%  '

# Row 15
$ws.Range('D15').Value = '3.664.54'
$ws.Range('E15').Value = '  +2.74%  '

# Row 16
$ws.Range('D16').Value = '65.011.85'
$ws.Range('E16').Value = '  +1.60%  '

# Row 17
$ws.Range('D17').Value = '3.185.11'
$ws.Range('E17').Value = '  +4.06%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.41%  '

# Row 19
$ws.Range('E19').Value = '  +1.08%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '505.90'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.65%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.77'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.76%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.716'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.52%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.49%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.17%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.06'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.54%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.55%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.86'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.26%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.12%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '27.56'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.72%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.78'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.47%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.06%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.19'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.26%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.19'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.46%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.51'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.20%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.10'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.51%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0894'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.67%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '464.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.41%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0420'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.84%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.99'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.12%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.64'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.45%  '

# Row 42
$ws.Range('D42').Value = '3.055.59'
$ws.Range('E42').Value = '  +0.62%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.117'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.04%  '

# Row 44
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.08%  '

# Row 45
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.281'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.91%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.60'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.84%  '

# Row 47
$ws.Range('D47').Value = '0.0₃0584'
$ws.Range('E47').Value = '  +12.45%  '

# Row 48
$ws.Range('E48').Value = '  +0.00%  '

# Row 49
$ws.Range('E49').Value = '  +0.44%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.09%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.31'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.20%  '
